# Update the two-digit-by-one-digit division problems in the single 20x5 table.
# Each of the 5 populated rows (1, 5, 9, 13, 17) has 5 cells; we target each cell
# explicitly by (row, column) rather than a document-wide Find/Replace, because
# several of the original expressions repeat verbatim across different cells and
# this host's Find.Execute matches the first occurrence in the whole document
# rather than honoring the scoping Range, which would corrupt unrelated cells.
$d = $word.ActiveDocument
$t = $d.Tables(1)

# 50÷9= -> 56÷5=
$t.Cell(1,1).Range.Text = "56÷5="
# 49÷6= -> 85÷9=
$t.Cell(1,2).Range.Text = "85÷9="
# 23÷4= -> 41÷6=
$t.Cell(1,3).Range.Text = "41÷6="
# 63÷2= -> 14÷2=
$t.Cell(1,4).Range.Text = "14÷2="
# 48÷6= -> 23÷2=
$t.Cell(1,5).Range.Text = "23÷2="
# 87÷3= -> 36÷8=
$t.Cell(5,1).Range.Text = "36÷8="
# 95÷3= -> 63÷8=
$t.Cell(5,2).Range.Text = "63÷8="
# 64÷2= -> 59÷7=
$t.Cell(5,3).Range.Text = "59÷7="
# 93÷8= -> 63÷4=
$t.Cell(5,4).Range.Text = "63÷4="
# 32÷9= -> 59÷4=
$t.Cell(5,5).Range.Text = "59÷4="
# 81÷3= -> 22÷9=
$t.Cell(9,1).Range.Text = "22÷9="
# 23÷9= -> 65÷8=
$t.Cell(9,2).Range.Text = "65÷8="
# 64÷2= -> 89÷5=
$t.Cell(9,3).Range.Text = "89÷5="
# 50÷7= -> 78÷3=
$t.Cell(9,4).Range.Text = "78÷3="
# 74÷5= -> 86÷6=
$t.Cell(9,5).Range.Text = "86÷6="
# 86÷8= -> 37÷3=
$t.Cell(13,1).Range.Text = "37÷3="
# 65÷5= -> 87÷4=
$t.Cell(13,2).Range.Text = "87÷4="
# 14÷2= -> 25÷3=
$t.Cell(13,3).Range.Text = "25÷3="
# 13÷2= -> 50÷4=
$t.Cell(13,4).Range.Text = "50÷4="
# 74÷3= -> 91÷8=
$t.Cell(13,5).Range.Text = "91÷8="
# 85÷9= -> 30÷3=
$t.Cell(17,1).Range.Text = "30÷3="
# 34÷8= -> 20÷2=
$t.Cell(17,2).Range.Text = "20÷2="
# 91÷2= -> 95÷9=
$t.Cell(17,3).Range.Text = "95÷9="
# 40÷9= -> 39÷5=
$t.Cell(17,4).Range.Text = "39÷5="
# 75÷9= -> 76÷7=
$t.Cell(17,5).Range.Text = "76÷7="
